# TriviaNight.pptx edit:
# Every "Question N of 10 in round M" / "Answer N of 10 in round M" text run
# had its question-number and round-number swapped (the round number keeps
# its original zero-padded width, the question number is left unpadded).
#
# Walk every shape on every slide and, whenever its text matches the
# "Question/Answer X of Y in round Z" pattern, swap X and Z.

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tf = $shp.TextFrame
            if ($tf.HasText) {
                $text = $tf.TextRange.Text
                if ($text -match '^(Question|Answer) (\d+) of (\d+) in round (\d+)$') {
                    $kind = $matches[1]
                    $questionNum = $matches[2]
                    $total = $matches[3]
                    $roundNum = $matches[4]

                    $roundWidth = $roundNum.Length
                    $newQuestionNum = [int]$roundNum
                    $newRoundNum = ([int]$questionNum).ToString().PadLeft($roundWidth, '0')

                    $newText = "$kind $newQuestionNum of $total in round $newRoundNum"

                    if ($newText -ne $text) {
                        $tf.TextRange.Text = $newText
                    }
                }
            }
        }
    }
}
